$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "readme" sheet / Table1: re-order the log columns so "Author" comes
#    right after "index" (before "sheet_name"), i.e. new header order is
#    index, Author, sheet_name, Date, JobNo.
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("readme")

# Header row (row 1): B=Author, C=sheet_name, D=Date (A=index, E=JobNo unchanged)
$readme.Cells.Item(1, 2).Value = "Author"
$readme.Cells.Item(1, 3).Value = "sheet_name"
$readme.Cells.Item(1, 4).Value = "Date"

# Data rows 2-12: rotate the B/C/D values so the data follows the header move
# (old B=sheet_name, old C=Date, old D=Author) -> (new B=Author, new C=sheet_name, new D=Date)
for ($r = 2; $r -le 12; $r++) {
    $oldB = $readme.Cells.Item($r, 2).Value()
    $oldC = $readme.Cells.Item($r, 3).Value()
    $oldD = $readme.Cells.Item($r, 4).Value()

    $readme.Cells.Item($r, 2).Value = $oldD
    $readme.Cells.Item($r, 3).Value = $oldB
    $readme.Cells.Item($r, 4).Value = $oldC
}

# ---------------------------------------------------------------------------
# 2) "Project Information" sheet: refresh the "Date of Analysis" timestamp.
# ---------------------------------------------------------------------------
$projInfo = $wb.Worksheets.Item("Project Information")
$projInfo.Cells.Item(11, 2).Value = "2022-03-08 14:52:54.463102"

# ---------------------------------------------------------------------------
# 3) Every "Results, Air Speed x" sheet/table: swap "Room Name"/"Room ID" so
#    Room ID becomes the first (index) column, Room Name the second.
# ---------------------------------------------------------------------------
$resultSheets = @(
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

foreach ($sheetName in $resultSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header row
    $oldA1 = $ws.Cells.Item(1, 1).Value()
    $oldB1 = $ws.Cells.Item(1, 2).Value()
    $ws.Cells.Item(1, 1).Value = $oldB1
    $ws.Cells.Item(1, 2).Value = $oldA1

    # Data rows 2-32
    for ($r = 2; $r -le 32; $r++) {
        $oldA = $ws.Cells.Item($r, 1).Value()
        $oldB = $ws.Cells.Item($r, 2).Value()
        $ws.Cells.Item($r, 1).Value = $oldB
        $ws.Cells.Item($r, 2).Value = $oldA
    }
}
